# TEST-cyclic.xlsx — "dynamic construction of ranges, index, match, indirect (wip)"
#
# Adds a small INDIRECT demo block to the "Sheet6" tab (physically
# xl/worksheets/sheet5.xml): a few label/value rows plus dynamic-array
# INDIRECT() lookups that resolve through a plain cell reference, a
# concatenated "B"&row reference, a cross-sheet "Sheet6!B"&row reference,
# and a workbook defined name ("George" -> Sheet6!$B$13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet6")

# --- workbook-level defined name: George -> Sheet6!$B$13 -------------------
# (single quotes so PowerShell doesn't try to expand "$B" as a variable)
$wb.Names.Add("George", '=Sheet6!$B$13')

# --- new label/value rows (11-14) ------------------------------------------
# Written in this order so the shared-strings table picks up the same
# indices as the authored file: George, then B11, then B12.
$ws.Range("A13").Value = "George"
$ws.Range("A11").Value = "B11"
$ws.Range("A12").Value = "B12"

$ws.Range("B11").Value = 1.333
$ws.Range("B12").Value = 45
$ws.Range("B13").Value = 10
$ws.Range("A14").Value = 14
$ws.Range("B14").Value = 62

# --- dynamic-array INDIRECT() formulas (16-20) ------------------------------
$ws.Range("A16").FormulaArray = "=INDIRECT(A11)"
$ws.Range("A17").FormulaArray = "=INDIRECT(A12)"
$ws.Range("A18").FormulaArray = "=INDIRECT(A13)"
$ws.Range("A19").FormulaArray = '=INDIRECT("B"&A14)'
$ws.Range("A20").FormulaArray = '=INDIRECT("Sheet6!B"&A14)'

# --- view state: selection moves to A21, window scrolls into place ---------
$ws.Activate() | Out-Null
$ws.Range("A21").Select() | Out-Null

$win = $wb.Windows.Item(1)
$win.Left = 3840
$win.Top = 0
